$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: update Correspond Handoff Datetime (D4) and Correspond Handback DateTime (G4)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D4").Value = "2016-01-20 07:39:12"
$wsZh.Range("G4").Value = "2016-01-20 07:39:57"

# "de-de" sheet: update Correspond Handoff Datetime (D4) and Correspond Handback DateTime (G4)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D4").Value = "2016-01-20 07:39:22"
$wsDe.Range("G4").Value = "2016-01-20 07:40:17"
